# Apply the MovementSpd unit-fix edit described in the commit message:
# Character - MovementSpd 5.75 -> 575 (cell F3 on sheet "CharacterGameData")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharacterGameData")

# Update MovementSpd value for the character row (F3)
$ws.Range("F3").Value = 575

# Reflect the selection shown in the saved file (cell I6 was last selected)
$ws.Activate()
$ws.Range("I6").Select()
